$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion rate text in cell A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.34 = 12937.1 pesos`n✅ 12937.1 pesos = 3.33 = 964.74 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$wsHoja1.Range("A1").Value = $newText

# --- Sheet "tasas": update the computed rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")

$wsTasas.Range("N10").Value = 298.989
$wsTasas.Range("O10").Value = 3868.05
$wsTasas.Range("N12").Value = 3888.9
